# This workbook holds one weekly "Caqui" price-sheet; the data rows (2-16)
# get re-ordered/refreshed to reflect the latest weekly snapshot - every
# row's full content (A:T) ends up identical to some row's content from
# before the edit, just under a different row number. Re-create that by
# reading the old grid into memory first, then writing it back out in the
# new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 16
$lastCol = 20   # column T

# Snapshot every source cell (row -> column -> value) before we overwrite
# anything, since several destination rows read from rows we're about to
# clobber.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# newRow -> oldRow the content now comes from.
$rowMap = @{
    2  = 13
    3  = 2
    4  = 15
    5  = 4
    6  = 12
    7  = 10
    8  = 16
    9  = 11
    10 = 14
    11 = 5
    12 = 7
    13 = 8
    14 = 6
    15 = 3
    16 = 9
}

foreach ($newRow in ($rowMap.Keys | Sort-Object)) {
    $oldRow = $rowMap[$newRow]
    $source = $snapshot[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $source[$c]
    }
}
